{"js": "// Replace each two-digit multiplication prompt with its new value.\n// Every \"AA\u00d7BB=\" string in this worksheet is unique, so a plain\n// search-and-replace (no wildcards) against the whole document body is\n// safe and cannot cross-contaminate other cells.\nconst pairs = [\n  [\"24\u00d757=\", \"83\u00d765=\"],\n  [\"13\u00d796=\", \"77\u00d795=\"],\n  [\"23\u00d738=\", \"63\u00d785=\"],\n  [\"40\u00d759=\", \"86\u00d776=\"],\n  [\"43\u00d796=\", \"27\u00d733=\"],\n  [\"33\u00d714=\", \"92\u00d719=\"],\n  [\"87\u00d796=\", \"39\u00d764=\"],\n  [\"94\u00d756=\", \"62\u00d752=\"],\n  [\"89\u00d789=\", \"93\u00d728=\"],\n  [\"97\u00d794=\", \"74\u00d754=\"],\n  [\"26\u00d753=\", \"33\u00d711=\"],\n  [\"47\u00d731=\", \"31\u00d724=\"],\n  [\"55\u00d731=\", \"27\u00d715=\"],\n  [\"19\u00d772=\", \"23\u00d749=\"],\n  [\"15\u00d777=\", \"15\u00d787=\"],\n  [\"71\u00d729=\", \"35\u00d773=\"],\n  [\"72\u00d732=\", \"21\u00d750=\"],\n  [\"62\u00d729=\", \"65\u00d764=\"],\n  [\"14\u00d790=\", \"60\u00d738=\"],\n  [\"48\u00d711=\", \"56\u00d770=\"],\n  [\"74\u00d783=\", \"66\u00d739=\"],\n  [\"67\u00d791=\", \"12\u00d750=\"],\n  [\"50\u00d749=\", \"67\u00d799=\"],\n  [\"44\u00d730=\", \"70\u00d760=\"],\n  [\"95\u00d716=\", \"90\u00d742=\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication prompt with its new value.\n# Every \"AA\u00d7BB=\" string in this worksheet is unique, so a plain\n# Find/Replace (no wildcards) on the whole document body is safe and\n# cannot cross-contaminate other cells.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"24\u00d757=\"; New = \"83\u00d765=\" },\n    @{ Old = \"13\u00d796=\"; New = \"77\u00d795=\" },\n    @{ Old = \"23\u00d738=\"; New = \"63\u00d785=\" },\n    @{ Old = \"40\u00d759=\"; New = \"86\u00d776=\" },\n    @{ Old = \"43\u00d796=\"; New = \"27\u00d733=\" },\n    @{ Old = \"33\u00d714=\"; New = \"92\u00d719=\" },\n    @{ Old = \"87\u00d796=\"; New = \"39\u00d764=\" },\n    @{ Old = \"94\u00d756=\"; New = \"62\u00d752=\" },\n    @{ Old = \"89\u00d789=\"; New = \"93\u00d728=\" },\n    @{ Old = \"97\u00d794=\"; New = \"74\u00d754=\" },\n    @{ Old = \"26\u00d753=\"; New = \"33\u00d711=\" },\n    @{ Old = \"47\u00d731=\"; New = \"31\u00d724=\" },\n    @{ Old = \"55\u00d731=\"; New = \"27\u00d715=\" },\n    @{ Old = \"19\u00d772=\"; New = \"23\u00d749=\" },\n    @{ Old = \"15\u00d777=\"; New = \"15\u00d787=\" },\n    @{ Old = \"71\u00d729=\"; New = \"35\u00d773=\" },\n    @{ Old = \"72\u00d732=\"; New = \"21\u00d750=\" },\n    @{ Old = \"62\u00d729=\"; New = \"65\u00d764=\" },\n    @{ Old = \"14\u00d790=\"; New = \"60\u00d738=\" },\n    @{ Old = \"48\u00d711=\"; New = \"56\u00d770=\" },\n    @{ Old = \"74\u00d783=\"; New = \"66\u00d739=\" },\n    @{ Old = \"67\u00d791=\"; New = \"12\u00d750=\" },\n    @{ Old = \"50\u00d749=\"; New = \"67\u00d799=\" },\n    @{ Old = \"44\u00d730=\"; New = \"70\u00d760=\" },\n    @{ Old = \"95\u00d716=\"; New = \"90\u00d742=\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.New\n    $find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n"}
